$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2022-08-02 14:45:05"
$newTimestamp = "2022-08-02 20:57:33"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 15).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 73 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 15)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value2 = $newTimestamp
    }
}
